function Find-ParaIndex($d, $pattern, $start) {
    for ($i = $start; $i -le $d.Paragraphs.Count; $i++) {
        if ($d.Paragraphs($i).Range.Text -like $pattern) {
            return $i
        }
    }
    return -1
}

function Set-RangeXml($range, $innerBodyXml) {
    $xml = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + $innerBodyXml + '</w:document></pkg:xmlData></pkg:part></pkg:package>'
    $range.InsertXML($xml)
}

$d = $word.ActiveDocument

# --- 1. Split the red title paragraph into an empty paragraph + the corrected text ---
$titleIdx = Find-ParaIndex $d "*Crear un men* responsive*" 1
$titleRange = $d.Paragraphs($titleIdx).Range
$titleBody = '<w:body>' +
  '<w:p><w:pPr><w:jc w:val="center"/><w:rPr><w:b/><w:bCs/><w:color w:val="FF0000"/></w:rPr></w:pPr></w:p>' +
  '<w:p><w:pPr><w:jc w:val="center"/><w:rPr><w:b/><w:bCs/><w:color w:val="FF0000"/></w:rPr></w:pPr>' +
  '<w:r><w:rPr><w:b/><w:bCs/><w:color w:val="FF0000"/></w:rPr><w:t xml:space="preserve">Crear un men&#250; responsive con </w:t></w:r>' +
  '<w:proofErr w:type="spellStart"/>' +
  '<w:r><w:rPr><w:b/><w:bCs/><w:color w:val="FF0000"/></w:rPr><w:t>Tailwind</w:t></w:r>' +
  '<w:proofErr w:type="spellEnd"/>' +
  '<w:r><w:rPr><w:b/><w:bCs/><w:color w:val="FF0000"/></w:rPr><w:t xml:space="preserve"> y Alpine para el blog</w:t></w:r>' +
  '</w:p>' +
  '</w:body>'
Set-RangeXml $titleRange $titleBody

# --- 2. Reword the "Alpine tratara el div..." bullet so "div" gets its own spell-checked run ---
$alpineIdx = Find-ParaIndex $d "*Alpine trata* el div respectivo*" 1
$alpineRange = $d.Paragraphs($alpineIdx).Range
$alpineBody = '<w:body><w:p><w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr></w:pPr>' +
  '<w:r><w:t xml:space="preserve"> Alpine tratar&#225; el </w:t></w:r>' +
  '<w:proofErr w:type="spellStart"/>' +
  '<w:r><w:t>div</w:t></w:r>' +
  '<w:proofErr w:type="spellEnd"/>' +
  '<w:r><w:t xml:space="preserve"> respectivo como un componente. Para lograrlo, es necesario agregarle al div la propiedad x-data, de la siguiente manera:</w:t></w:r>' +
  '</w:p></w:body>'
Set-RangeXml $alpineRange $alpineBody

# --- 3. Mark the last five inline images as not-proofed (adds <w:rPr><w:noProof/></w:rPr> to their runs) ---
for ($i = $d.InlineShapes.Count - 4; $i -le $d.InlineShapes.Count; $i++) {
    $shapeRange = $d.InlineShapes.Item($i).Range
    $shapeRange.NoProofing = 1
}

# --- 4. Drop the leftover _GoBack bookmark paragraph ---
$gobackIdx = Find-ParaIndex $d "*cierre:*" 1
$gobackIdx = $gobackIdx + 1
$gobackRange = $d.Paragraphs($gobackIdx).Range
$gobackBody = '<w:body><w:p><w:pPr><w:ind w:left="360"/><w:jc w:val="both"/></w:pPr></w:p></w:body>'
Set-RangeXml $gobackRange $gobackBody

Write-Host "Done"
